$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.772.86'
$ws.Range("D3").Value = '3.417.39'
$ws.Range("E3").Value = '  -3.56%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'580.37"
$ws.Range("E5").Value = '  -4.13%  '
$ws.Range("D6").Value = "'133.67"
$ws.Range("E6").Value = '  -7.66%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.417.11'
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("E9").Value = '  -6.52%  '
$ws.Range("D10").Value = "'0.120"
$ws.Range("E10").Value = '  -8.80%  '
$ws.Range("D11").Value = "'7.03"
$ws.Range("E11").Value = '  -10.33%  '
$ws.Range("D12").Value = "'0.372"
$ws.Range("E12").Value = '  -9.49%  '
$ws.Range("D13").Value = '3.995.83'
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("E14").Value = '  -9.30%  '
$ws.Range("D16").Value = '3.414.11'
$ws.Range("E16").Value = '  -3.62%  '
$ws.Range("E17").Value = '  -9.34%  '
$ws.Range("D18").Value = '64.676.05'
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").Value = "'9.42"
$ws.Range("E19").Value = '  -15.55%  '
$ws.Range("E20").Value = '  -8.01%  '
$ws.Range("D21").Value = "'13.39"
$ws.Range("E21").Value = '  -8.91%  '
$ws.Range("D22").Value = "'379.32"
$ws.Range("E22").Value = '  -10.31%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = "'0.538"
$ws.Range("E24").Value = '  -9.85%  '
$ws.Range("D25").Value = "'71.63"
$ws.Range("E25").Value = '  -7.41%  '
$ws.Range("D26").Value = '3.553.39'
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("E27").Value = '  -10.50%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("D29").Value = "'7.12"
$ws.Range("E29").Value = '  -9.98%  '
$ws.Range("E30").Value = '  -11.94%  '
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = '  -11.63%  '
$ws.Range("D32").Value = '3.436.42'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = "'0.141"
$ws.Range("E34").Value = '  -9.30%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'22.77"
$ws.Range("E35").Value = '  -6.25%  '
$ws.Range("D36").Value = "'169.19"
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = '  -13.48%  '
$ws.Range("E38").Value = '  -13.37%  '
$ws.Range("D39").Value = "'1.44"
$ws.Range("E39").Value = '  -12.51%  '
$ws.Range("E40").Value = '  -14.74%  '
$ws.Range("D41").Value = "'0.0754"
$ws.Range("E41").Value = '  -8.32%  '
$ws.Range("D42").Value = "'0.801"
$ws.Range("E42").Value = '  -6.96%  '
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = "'41.86"
$ws.Range("E44").Value = '  -7.93%  '
$ws.Range("E45").Value = '  -15.43%  '
$ws.Range("E46").Value = '  -9.81%  '
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = "'22.19"
$ws.Range("E48").Value = '  -6.97%  '
$ws.Range("D49").Value = "'6.42"
$ws.Range("E49").Value = '  -9.31%  '
$ws.Range("D50").Value = '2.186.85'
$ws.Range("E50").Value = '  -6.22%  '
$ws.Range("D51").Value = "'1.95"
$ws.Range("E51").Value = '  -18.73%  '

# Reset style on forced-text numeric-looking cells to avoid stray quote-prefix style
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
